$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.028.84'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.754.80'
$ws.Range('E3').Value = '  -2.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.85'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3770'
$ws.Range('E7').Value = '  -4.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3346'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.34'
$ws.Range('E9').Value = '  -6.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.115'
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07192'
$ws.Range('E11').Value = '  -4.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.54'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.133'
$ws.Range('E14').Value = '  -5.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.142'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = '1.758.42'
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001054'
$ws.Range('E17').Value = '  -4.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06593'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '80.53'
$ws.Range('E19').Value = '  -5.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.89'
$ws.Range('E21').Value = '  -4.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.235'
$ws.Range('E22').Value = '  -4.93%  '
$ws.Range('D23').Value = '28.053.65'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.65'
$ws.Range('E24').Value = '  -5.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.398'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.84'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.85'
$ws.Range('E27').Value = '  -7.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.320'
$ws.Range('E28').Value = '  -7.70%  '
$ws.Range('D29').Value = '1.957.57'
$ws.Range('E29').Value = '  -2.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.77'
$ws.Range('E30').Value = '  -2.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.251'
$ws.Range('E31').Value = '  -15.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.021'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.769'
$ws.Range('E33').Value = '  -6.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08742'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.19'
$ws.Range('E35').Value = '  -6.99%  '
$ws.Range('B36').Value = 'TheSandbox'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6678'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02330'
$ws.Range('E37').Value = '  -3.93%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.163'
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06170'
$ws.Range('E39').Value = '  -5.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2106'
$ws.Range('E40').Value = '  -4.53%  '
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.447'
$ws.Range('E42').Value = '  -10.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.017'
$ws.Range('E43').Value = '  -5.50%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.67'
$ws.Range('E45').Value = '  -6.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.833'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6032'
$ws.Range('E47').Value = '  -5.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '128.62'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.016'
$ws.Range('E49').Value = '  -5.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.180'
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07169'
$ws.Range('E51').Value = '  -0.27%  '
